$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial value that was bumped by one
# day (2023-10-06 -> 2023-10-07, i.e. serial 45205 -> 45206) for every
# data row (rows 2 through 206).
$ws.Range("C2:C206").Value = 45206
